$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 11 (Item ID 5533)
$ws.Range("H11").Value = 116858.4
$ws.Range("I11").Value = 116858.4
$ws.Range("K11").Value = 116858.4
$ws.Range("M11").Value = -116718.4

# Row 74 (Item ID 5507)
$ws.Range("H74").Value = 11822.056
$ws.Range("I74").Value = 13319.934
$ws.Range("J74").Value = 4332.6665
$ws.Range("K74").Value = 13319.934
$ws.Range("L74").Value = 4332.6665
$ws.Range("M74").Value = -12383.934
$ws.Range("N74").Value = -6204.6665

# Row 77 (Item ID 5507)
$ws.Range("H77").Value = 11822.056
$ws.Range("I77").Value = 13319.934
$ws.Range("J77").Value = 4332.6665
$ws.Range("K77").Value = 66599.67
$ws.Range("L77").Value = 21663.3325
$ws.Range("M77").Value = -61919.67
$ws.Range("N77").Value = -31023.3325

# Row 86 (Item ID 12603)
$ws.Range("H86").Value = 95240500
$ws.Range("I86").Value = 90911416
$ws.Range("J86").Value = 111113784
$ws.Range("K86").Value = 90911416
$ws.Range("L86").Value = 111113784
$ws.Range("M86").Value = -90910293
$ws.Range("N86").Value = -111116030

# Row 88 (Item ID 12608)
$ws.Range("H88").Value = 1462.5
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()

# Row 89 (Item ID 12603)
$ws.Range("H89").Value = 95240500
$ws.Range("I89").Value = 90911416
$ws.Range("J89").Value = 111113784
$ws.Range("K89").Value = 454557080
$ws.Range("L89").Value = 555568920
$ws.Range("M89").Value = -454551464
$ws.Range("N89").Value = -555580152

# Row 91 (Item ID 12608)
$ws.Range("H91").Value = 1462.5
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()

# Row 137 (Item ID 44013)
$ws.Range("H137").Value = 4547496
$ws.Range("I137").Value = 1512.3334
$ws.Range("K137").Value = 4537.0002
$ws.Range("M137").Value = -1987.0002

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (Item ID 44147)
$ws.Range("H32").Value = 277396.78
$ws.Range("I32").Value = 298049.2
$ws.Range("K32").Value = 298049.2
$ws.Range("M32").Value = -297762.2

# Row 63 (Item ID 12528)
$ws.Range("H63").Value = 15910.692
$ws.Range("I63").Value = 7312.1665
$ws.Range("J63").Value = 23280.857
$ws.Range("K63").Value = 7312.1665
$ws.Range("L63").Value = 23280.857
$ws.Range("M63").Value = -6626.1665
$ws.Range("N63").Value = -24652.857

# Row 66 (Item ID 12528)
$ws.Range("H66").Value = 15910.692
$ws.Range("I66").Value = 7312.1665
$ws.Range("J66").Value = 23280.857
$ws.Range("K66").Value = 36560.8325
$ws.Range("L66").Value = 116404.285
$ws.Range("M66").Value = -33128.8325
$ws.Range("N66").Value = -123268.285

# Row 88 (Item ID 12530)
$ws.Range("H88").Value = 2462.9167
$ws.Range("J88").Value = 3482.3333
$ws.Range("L88").Value = 3482.3333
$ws.Range("N88").Value = -4294.3333

# Row 91 (Item ID 12530)
$ws.Range("H91").Value = 2462.9167
$ws.Range("J91").Value = 3482.3333
$ws.Range("L91").Value = 3482.3333
$ws.Range("N91").Value = -6290.3333

# Row 132 (Item ID 43997)
$ws.Range("H132").Value = 3171.7026
$ws.Range("I132").Value = 1137.5
$ws.Range("K132").Value = 3412.5
$ws.Range("M132").Value = -882.5

$ws = $wb.Worksheets.Item("BSM")
# Row 20 (Item ID 14149)
$ws.Range("H20").Value = 1405.3684
$ws.Range("I20").Value = 1329.5454
$ws.Range("J20").Value = 1509.625
$ws.Range("K20").Value = 1329.5454
$ws.Range("L20").Value = 1509.625
$ws.Range("M20").Value = -1082.5454
$ws.Range("N20").Value = -2003.625

# Row 82 (Item ID 11877)
$ws.Range("H82").Value = 8032.8335
$ws.Range("I82").Value = 8032.8335
$ws.Range("K82").Value = 8032.8335
$ws.Range("M82").Value = -7649.8335

# Row 85 (Item ID 11877)
$ws.Range("H85").Value = 8032.8335
$ws.Range("I85").Value = 8032.8335
$ws.Range("K85").Value = 8032.8335
$ws.Range("M85").Value = -6706.8335

# Row 86 (Item ID 12526)
$ws.Range("H86").Value = 3195.4783
$ws.Range("I86").Value = 2084.1538
$ws.Range("J86").Value = 4640.2
$ws.Range("K86").Value = 2084.1538
$ws.Range("L86").Value = 4640.2
$ws.Range("M86").Value = -961.1538
$ws.Range("N86").Value = -6886.2

# Row 89 (Item ID 12526)
$ws.Range("H89").Value = 3195.4783
$ws.Range("I89").Value = 2084.1538
$ws.Range("J89").Value = 4640.2
$ws.Range("K89").Value = 10420.769
$ws.Range("L89").Value = 23201
$ws.Range("M89").Value = -4804.769
$ws.Range("N89").Value = -34433

# Row 134 (Item ID 43998)
$ws.Range("H134").Value = 30003128
$ws.Range("I134").Value = 2492.5908
$ws.Range("K134").Value = 7477.7724
$ws.Range("M134").Value = -4942.7724

$ws = $wb.Worksheets.Item("CRP")
# Row 58 (Item ID 44021)
$ws.Range("H58").Value = 4379.3335
$ws.Range("J58").Value = 4564.25
$ws.Range("L58").Value = 4564.25
$ws.Range("N58").Value = -4970.25

# Row 136 (Item ID 44021)
$ws.Range("H136").Value = 4379.3335
$ws.Range("J136").Value = 4564.25
$ws.Range("L136").Value = 13692.75
$ws.Range("N136").Value = -18792.75

$ws = $wb.Worksheets.Item("CUL")
# Row 113 (Item ID 27843)
$ws.Range("H113").Value = 1185.9
$ws.Range("J113").Value = 1138.4286
$ws.Range("L113").Value = 3415.2858
$ws.Range("N113").Value = -7755.2858

# Row 129 (Item ID 36054)
$ws.Range("H129").Value = 19755244
$ws.Range("I129").Value = 1318.3334
$ws.Range("K129").Value = 3955.0002
$ws.Range("M129").Value = 1044.9998

# Row 131 (Item ID 36060)
$ws.Range("H131").Value = 4514259.5
$ws.Range("I131").Value = 9093696
$ws.Range("J131").Value = 3205849.2
$ws.Range("K131").Value = 27281088
$ws.Range("L131").Value = 9617547.600000001
$ws.Range("M131").Value = -27276048
$ws.Range("N131").Value = -9627627.600000001

$ws = $wb.Worksheets.Item("GSM")
# Row 70 (Item ID 14146)
$ws.Range("H70").Value = 9824.257
$ws.Range("I70").Value = 9395.482
$ws.Range("K70").Value = 9395.482
$ws.Range("M70").Value = -9125.482

# Row 73 (Item ID 14146)
$ws.Range("H73").Value = 9824.257
$ws.Range("I73").Value = 9395.482
$ws.Range("K73").Value = 9395.482
$ws.Range("M73").Value = -8459.482

# Row 100 (Item ID 18367)
$ws.Range("H100").Value = 84000
$ws.Range("J100").Value = 84000
$ws.Range("L100").Value = 84000
$ws.Range("N100").Value = -86164

# Row 102 (Item ID 36169)
$ws.Range("H102").Value = 21740084
$ws.Range("I102").Value = 23810506
$ws.Range("K102").Value = 23810506
$ws.Range("M102").Value = -23808884

$ws = $wb.Worksheets.Item("LTW")
# Row 7 (Item ID 36249)
$ws.Range("H7").Value = 11550.0625
$ws.Range("I7").Value = 4641
$ws.Range("J7").Value = 15695.5
$ws.Range("K7").Value = 4641
$ws.Range("L7").Value = 15695.5
$ws.Range("M7").Value = -4529
$ws.Range("N7").Value = -15919.5

# Row 40 (Item ID 36248)
$ws.Range("H40").Value = 5220.7856
$ws.Range("I40").Value = 4475
$ws.Range("K40").Value = 4475
$ws.Range("M40").Value = -4339

# Row 43 (Item ID 4314)
$ws.Range("H43").Value = 1073466.1
$ws.Range("J43").Value = 1248159.2
$ws.Range("L43").Value = 1248159.2
$ws.Range("N43").Value = -1248545.2

# Row 55 (Item ID 5284)
$ws.Range("H55").Value = 212.52942
$ws.Range("J55").Value = 234.92308
$ws.Range("L55").Value = 234.92308
$ws.Range("N55").Value = -580.92308

# Row 62 (Item ID 10740)
$ws.Range("H62").Value = 46500
$ws.Range("J62").Value = 46500
$ws.Range("L62").Value = 46500
$ws.Range("N62").Value = -47748

# Row 65 (Item ID 10740)
$ws.Range("H65").Value = 46500
$ws.Range("J65").Value = 46500
$ws.Range("L65").Value = 139500
$ws.Range("N65").Value = -145740

# Row 82 (Item ID 12565)
$ws.Range("H82").Value = 3513.9048
$ws.Range("I82").Value = 1157.3077
$ws.Range("K82").Value = 1157.3077
$ws.Range("M82").Value = -796.3077000000001

# Row 85 (Item ID 12565)
$ws.Range("H85").Value = 3513.9048
$ws.Range("I85").Value = 1157.3077
$ws.Range("K85").Value = 1157.3077
$ws.Range("M85").Value = 90.69229999999993

# Row 114 (Item ID 25990)
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").ClearContents()
$ws.Range("N114").Value = 0

# Row 122 (Item ID 36247)
$ws.Range("H122").Value = 4174.34
$ws.Range("I122").Value = 3099
$ws.Range("K122").Value = 9297
$ws.Range("M122").Value = -6847

# Row 126 (Item ID 36249)
$ws.Range("H126").Value = 11550.0625
$ws.Range("I126").Value = 4641
$ws.Range("J126").Value = 15695.5
$ws.Range("K126").Value = 13923
$ws.Range("L126").Value = 47086.5
$ws.Range("M126").Value = -11453
$ws.Range("N126").Value = -52026.5

# Row 136 (Item ID 44060)
$ws.Range("H136").Value = 2088.7856
$ws.Range("I136").Value = 1875.25
$ws.Range("J136").Value = 3370
$ws.Range("K136").Value = 5625.75
$ws.Range("L136").Value = 10110
$ws.Range("M136").Value = -3075.75
$ws.Range("N136").Value = -15210

$ws = $wb.Worksheets.Item("WVR")
# Row 20 (Item ID 3023)
$ws.Range("H20").Value = 60
$ws.Range("I20").Value = 60
$ws.Range("K20").Value = 60
$ws.Range("M20").Value = 180

# Row 81 (Item ID 12596)
$ws.Range("H81").Value = 5900
$ws.Range("I81").Value = 5900
$ws.Range("K81").Value = 11800
$ws.Range("M81").Value = -10739

# Row 84 (Item ID 12596)
$ws.Range("H84").Value = 5900
$ws.Range("I84").Value = 5900
$ws.Range("K84").Value = 59000
$ws.Range("M84").Value = -53696

# Row 107 (Item ID 27746)
$ws.Range("H107").Value = 58824200
$ws.Range("I107").Value = 738.2308
$ws.Range("K107").Value = 2214.6924
$ws.Range("M107").Value = -294.6923999999999

# Row 122 (Item ID 36208)
$ws.Range("H122").Value = 3496.5938
$ws.Range("I122").Value = 3594.7778
$ws.Range("K122").Value = 10784.3334
$ws.Range("M122").Value = -8334.3334

# Row 136 (Item ID 44031)
$ws.Range("H136").Value = 37529.934
$ws.Range("J136").Value = 4900.45
$ws.Range("L136").Value = 14701.35
$ws.Range("N136").Value = -19801.35
